# Actualizacion automatica del tracker
# - Rellena el resultado del pick de la fila 213 (Zdenek Kolar vs Murkel Dellien)
# - Agrega seis picks nuevos del 2025-08-16 (filas 216-221)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resultado fila 213 ---
$ws.Range("G213").Value = "Fallo"
$ws.Range("H213").Value = -1

# --- Nuevas filas 216-221 ---

$rows = @(
    @{ Row=216; A=14453337; C="Amandine Hesse";      D="Diane Parry";      E="Gana Amandine Hesse";      F=4 },
    @{ Row=217; A=14452016; C="Emily Appleton";       D="Anna Rogers";      E="Gana Anna Rogers";         F=2 },
    @{ Row=218; A=14452021; C="Talia Gibson";         D="Vivian Wolff";     E="Gana Vivian Wolff";        F=3.75 },
    @{ Row=219; A=14452017; C="Aleksandra Krunic";    D="Yafan Wang";       E="Gana Aleksandra Krunic";   F=2.5 },
    @{ Row=220; A=14369089; C="Arthur Fery";          D="Bernard Tomic";    E="Gana Bernard Tomic";       F=2.63 },
    @{ Row=221; A=14370751; C="Jaime Faria";          D="Mattia Bellucci";  E="Gana Jaime Faria";         F=2.63 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.A
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
}

# La columna "fecha" (B) se guarda como texto en el tracker (p.ej. "2025-08-16"),
# asi que hay que evitar que Excel la convierta automaticamente en una fecha.
$fechaRange = $ws.Range("B216:B221")
$fechaRange.NumberFormat = "@"
$fechaRange.Value = "2025-08-16"
$fechaRange.Style = "Normal"

"Tracker actualizado"
